$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "plots"
$ws.Range("C1").Value = "block"
$ws.Range("D1").Value = "imb"
$ws.Range("E1").Value = "imbt"
